$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first data row (old row 2, date 2007-11-14 / 39400) was removed from the
# naive forecaster output. Deleting it shifts every subsequent row up by one,
# which reproduces the corrected date/year (A, B, D) columns and also fixes
# the sheet dimension to A1:E52.
$ws.Rows(2).Delete()

# The y_1 (column C) AR(2) fitted values were recomputed by the bugfixed
# naive forecaster component; write the corrected series back into C2:C52.
$cValues = New-Object 'object[,]' 51,1
$cValues[0,0] = $null
$cValues[1,0] = $null
$cValues[2,0] = $null
$cValues[3,0] = 0.1715429114845124
$cValues[4,0] = $null
$cValues[5,0] = 0.8004663283405655
$cValues[6,0] = $null
$cValues[7,0] = 5.253783907501819
$cValues[8,0] = $null
$cValues[9,0] = 3.522405026196918
$cValues[10,0] = 1.133560223479058
$cValues[11,0] = 1.656063945467268
$cValues[12,0] = 3.633318781899142
$cValues[13,0] = 4.06235252733802
$cValues[14,0] = 3.057638025163611
$cValues[15,0] = 3.05427116350534
$cValues[16,0] = 2.319057151538662
$cValues[17,0] = 2.305809238174006
$cValues[18,0] = 2.536029549059826
$cValues[19,0] = 2.509111342826809
$cValues[20,0] = 3.025024236774643
$cValues[21,0] = 3.120740332206995
$cValues[22,0] = 3.279355759764568
$cValues[23,0] = 3.296731496509198
$cValues[24,0] = 2.945303709067959
$cValues[25,0] = 2.891533899000343
$cValues[26,0] = 2.827707622797226
$cValues[27,0] = 2.861315725866587
$cValues[28,0] = 2.552688975800033
$cValues[29,0] = 2.618329006605924
$cValues[30,0] = 1.790319754067715
$cValues[31,0] = 1.790319754067715
$cValues[32,0] = 2.130407351599706
$cValues[33,0] = 2.137626121054947
$cValues[34,0] = 2.339531676162721
$cValues[35,0] = 2.339531676162721
$cValues[36,0] = 5.037171918133976
$cValues[37,0] = 4.951039758187648
$cValues[38,0] = 4.834496776263886
$cValues[39,0] = 4.834496776263886
$cValues[40,0] = 3.153537734543965
$cValues[41,0] = 2.838865660558509
$cValues[42,0] = 2.798216547494237
$cValues[43,0] = 2.798216547494237
$cValues[44,0] = 1.831762447564067
$cValues[45,0] = 1.625773169906108
$cValues[46,0] = 1.530879676868468
$cValues[47,0] = 1.530879676868468
$cValues[48,0] = 2.104676416355189
$cValues[49,0] = 2.030491763452114
$cValues[50,0] = 2.060859685319461
$ws.Range("C2:C52").Value = $cValues

# The y_1_forecast (column E) values were likewise recomputed; write the
# corrected series back into E2:E52.
$eValues = New-Object 'object[,]' 51,1
$eValues[0,0] = $null
$eValues[1,0] = $null
$eValues[2,0] = $null
$eValues[3,0] = $null
$eValues[4,0] = $null
$eValues[5,0] = $null
$eValues[6,0] = $null
$eValues[7,0] = $null
$eValues[8,0] = $null
$eValues[9,0] = 0.5784444854042281
$eValues[10,0] = 1.985690391709771
$eValues[11,0] = 2.529895848567842
$eValues[12,0] = 2.715291551682419
$eValues[13,0] = 4.060884847379076
$eValues[14,0] = 2.42782168586293
$eValues[15,0] = 2.270469368501771
$eValues[16,0] = 2.508920621023392
$eValues[17,0] = 2.467161166346266
$eValues[18,0] = 2.546671316138061
$eValues[19,0] = 2.480855794925163
$eValues[20,0] = 2.69389938681992
$eValues[21,0] = 2.775533179497169
$eValues[22,0] = 3.107596903291299
$eValues[23,0] = 3.221757900820066
$eValues[24,0] = 2.591074440292807
$eValues[25,0] = 2.545843589346886
$eValues[26,0] = 2.413544192054795
$eValues[27,0] = 2.631992339577627
$eValues[28,0] = 2.616345720823721
$eValues[29,0] = 2.671430903007876
$eValues[30,0] = 1.691013991470625
$eValues[31,0] = 2.153309886824961
$eValues[32,0] = 2.785334366326175
$eValues[33,0] = 2.891950990452763
$eValues[34,0] = 3.941556826710224
$eValues[35,0] = 4.667362054855917
$eValues[36,0] = 3.641364543513781
$eValues[37,0] = 3.481452844954491
$eValues[38,0] = 2.845322256798233
$eValues[39,0] = 3.305715257492858
$eValues[40,0] = 2.935215611250452
$eValues[41,0] = 2.377254777217375
$eValues[42,0] = 2.138412043368865
$eValues[43,0] = 1.757655717321982
$eValues[44,0] = 2.69124964061378
$eValues[45,0] = 2.42082970885531
$eValues[46,0] = 2.01742511619909
$eValues[47,0] = 2.159361127638926
$eValues[48,0] = 2.754798876280251
$eValues[49,0] = 2.559374235215039
$eValues[50,0] = 2.733459627814305
$ws.Range("E2:E52").Value = $eValues
